$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 199, shifting the existing
# rows 199-234 down to 200-235 (preserving all of their data/formatting).
$ws.Rows.Item(199).Insert()

# Populate the newly inserted row 199 with the new weekly data point.
$ws.Cells.Item(199, 1).Value = 5
$ws.Cells.Item(199, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(199, 3).Value = "Maule"
$ws.Cells.Item(199, 4).Value = 44522
$ws.Cells.Item(199, 5).Value = 7
$ws.Cells.Item(199, 6).Value = 100112032
$ws.Cells.Item(199, 7).Value = "Zapallo italiano"
$ws.Cells.Item(199, 8).Value = "Sin especificar"
$ws.Cells.Item(199, 9).Value = "Primera"
$ws.Cells.Item(199, 10).Value = 500
$ws.Cells.Item(199, 11).Value = 4000
$ws.Cells.Item(199, 12).Value = 4000
$ws.Cells.Item(199, 13).Value = 4000
$ws.Cells.Item(199, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(199, 15).Value = "Región del Maule"
$ws.Cells.Item(199, 16).Value = 67
$ws.Cells.Item(199, 17).Value = 60
$ws.Cells.Item(199, 18).Value = "Hortaliza"
